# Apply the "Updated cryptos list" data refresh (Thu Jun 15 14:24:46 UTC 2023).
# Row 42 (PaxosStandard) was delisted from the source feed: every row from 43
# downward shifted up by one, and a new "USDD" row was appended at the end (51).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new literal text value.
$updates = [ordered]@{
    'D2' = '25.010.56'
    'E2' = '  -3.83%  '
    'D3' = '1.643.97'
    'E3' = '  -5.65%  '
    'D4' = '0.9980'
    'E4' = '  -0.18%  '
    'D5' = '233.94'
    'E5' = '  -5.26%  '
    'D6' = '0.9995'
    'E6' = '  -0.07%  '
    'D7' = '0.4785'
    'E7' = '  -5.36%  '
    'E8' = '  -5.37%  '
    'D9' = '0.06122'
    'E9' = '  -1.00%  '
    'D10' = '0.07078'
    'E10' = '  -2.38%  '
    'D11' = '1.644.01'
    'E11' = '  -5.69%  '
    'D12' = '14.65'
    'E12' = '  -3.34%  '
    'D13' = '0.5970'
    'E13' = '  -8.69%  '
    'D14' = '4.357'
    'E14' = '  -6.97%  '
    'D15' = '73.73'
    'E15' = '  -5.01%  '
    'D16' = '0.9994'
    'E16' = '  -0.07%  '
    'D17' = '0.9993'
    'E17' = '  -0.02%  '
    'D18' = '25.009.21'
    'E18' = '  -3.89%  '
    'D19' = '0.000006600'
    'E19' = '  -3.98%  '
    'D20' = '11.27'
    'E20' = '  -5.36%  '
    'D21' = '1.852.69'
    'E21' = '  -6.13%  '
    'D22' = '4.346'
    'E22' = '  -2.94%  '
    'D23' = '8.590'
    'E23' = '  -1.52%  '
    'D24' = '5.248'
    'E24' = '  -2.37%  '
    'D25' = '134.40'
    'E25' = '  -1.00%  '
    'E26' = '  -2.46%  '
    'D27' = '1.387'
    'E27' = '  -7.75%  '
    'D28' = '103.95'
    'E28' = '  -1.46%  '
    'D29' = '1.651'
    'E29' = '  -7.54%  '
    'D30' = '3.883'
    'E30' = '  -0.65%  '
    'D31' = '0.07691'
    'E31' = '  -5.84%  '
    'D32' = '3.563'
    'E32' = '  -2.91%  '
    'D33' = '0.9988'
    'E33' = '  -0.02%  '
    'D34' = '0.04295'
    'E34' = '  -8.32%  '
    'D35' = '2.569'
    'E35' = '  -3.30%  '
    'E36' = '  -6.91%  '
    'D37' = '0.5917'
    'E37' = '  -3.48%  '
    'D38' = '2.578'
    'E38' = '  -6.68%  '
    'D39' = '0.8611'
    'E39' = '  +12.77%  '
    'E40' = '  -0.12%  '
    'D41' = '0.01513'
    'E41' = '  -6.90%  '
    'B42' = 'Quant'
    'C42' = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
    'D42' = '98.60'
    'E42' = '  -2.25%  '
    'B43' = 'RenderToken'
    'C43' = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
    'D43' = '1.768'
    'E43' = '  -8.62%  '
    'B44' = 'TheSandbox'
    'C44' = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
    'D44' = '0.3707'
    'E44' = '  -5.71%  '
    'B45' = 'FraxShare'
    'C45' = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
    'D45' = '4.685'
    'E45' = '  -6.61%  '
    'B46' = 'Algorand'
    'C46' = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
    'D46' = '0.1097'
    'E46' = '  -5.48%  '
    'B47' = 'Aptos'
    'C47' = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
    'D47' = '6.094'
    'E47' = '  -3.84%  '
    'B48' = 'Cronos'
    'C48' = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
    'D48' = '0.05203'
    'E48' = '  -1.82%  '
    'B49' = 'Elrond'
    'C49' = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
    'D49' = '28.98'
    'E49' = '  -5.60%  '
    'B50' = 'TrueUSD'
    'C50' = 'https://coinranking.com/coin/1ZZI6g5k5royD+trueusd-tusd'
    'D50' = '0.9993'
    'E50' = '  -0.27%  '
    'B51' = 'USDD'
    'C51' = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
    'D51' = '0.9979'
    'E51' = '  +0.01%  '
}

# Cells whose new text is purely numeric-looking (e.g. "233.94", "0.9980",
# "0.000006600") must be force-typed as Text first, otherwise Excel's COM
# layer auto-converts the assignment to a Number and silently drops
# formatting such as trailing zeros - these columns are text in the sheet.
$forceText = @(
    'D4'
    'D5'
    'D6'
    'D7'
    'D9'
    'D10'
    'D12'
    'D13'
    'D14'
    'D15'
    'D16'
    'D17'
    'D19'
    'D20'
    'D22'
    'D23'
    'D24'
    'D25'
    'D27'
    'D28'
    'D29'
    'D30'
    'D31'
    'D32'
    'D33'
    'D34'
    'D35'
    'D37'
    'D38'
    'D39'
    'D41'
    'D42'
    'D43'
    'D44'
    'D45'
    'D46'
    'D47'
    'D48'
    'D49'
    'D50'
    'D51'
)

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    if ($forceText -contains $addr) {
        $cell.NumberFormat = '@'
        $cell.Value = $updates[$addr]
        $cell.Style = 'Normal'
    } else {
        $cell.Value = $updates[$addr]
    }
}
